# Comentando e apagando códigos inúteis
#
# Populate row 1 (A1:D1) of the active sheet with the sample data that was
# previously left empty, extending the used range from A1:A1 to A1:D1.
#   A1 -> "s"        (text)
#   B1 -> "Lâmpada"   (text)
#   C1 -> 0           (number)
#   D1 -> FALSE       (boolean)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "s"
$ws.Range("B1").Value = "Lâmpada"
$ws.Range("C1").Value = 0
$ws.Range("D1").Value = $false

# Keep the newly written cells on the default ("Normal") style, i.e. do not
# let them inherit the explicit column formatting (styles 1/2) that was set
# up on the empty columns.
$ws.Range("A1:D1").Style = "Normal"
